$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Mark tasks 23 ("Create a smartphone UI", row 22) and 42 ("publish on app store", row 33)
# as Done in the "Done" column (B). The dependent "Status" column (C) recalculates
# automatically from "Ready"/"Blocked" to "Done" via the worksheet's existing formula.
$ws.Range("B22").Value = "Y"
$ws.Range("B33").Value = "Y"

# Update the view state: scroll the window so row 2 is at the top and select B2.
try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # view-scroll state isn't critical to the data change; ignore if unsupported
}
$ws.Range("B2").Select()
